$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update odds in row 4 (Negeri Sembilan vs Johor DT)
# ---------------------------------------------------------------------------
$ws.Range("G4").Value  = 32
$ws.Range("H4").Value  = 9.75
$ws.Range("J4").Value  = 21
$ws.Range("K4").Value  = 4.1
$ws.Range("L4").Value  = 1.19
$ws.Range("R4").Value  = 4.15
$ws.Range("U4").Value  = 2.64
$ws.Range("V4").Value  = 1.46
$ws.Range("W4").Value  = 200
$ws.Range("X4").Value  = 350
$ws.Range("Y4").Value  = 175
$ws.Range("AA4").Value = 600
$ws.Range("AB4").Value = 450
$ws.Range("AC4").Value = 32
$ws.Range("AD4").Value = 35
$ws.Range("AE4").Value = 70
$ws.Range("AH4").Value = 15.5
$ws.Range("AI4").Value = 8
$ws.Range("AJ4").Value = 17.5
$ws.Range("AK4").Value = 5.9
$ws.Range("AL4").Value = 12.5
$ws.Range("AM4").Value = 50
$ws.Range("AN4").Value = 37
$ws.Range("AO4").Value = 350
$ws.Range("AT4").Value = 6.2
$ws.Range("AU4").Value = 15.5
$ws.Range("AW4").Value = 3.55
$ws.Range("AY4").Value = 16.5
$ws.Range("AZ4").Value = 5.5

# ---------------------------------------------------------------------------
# 2) Update odds in row 5 (Terengganu vs Penang)
# ---------------------------------------------------------------------------
$ws.Range("P5").Value = 4.05
$ws.Range("S5").Value = 1.25
$ws.Range("T5").Value = 3.6
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.93

# ---------------------------------------------------------------------------
# 3) Insert a new row above the current row 8 (Hougang vs Tanjong Pagar moves
#    down to row 9) and populate the new row 8 with the Damac vs Al Riyadh
#    match.
# ---------------------------------------------------------------------------
$ws.Rows("8:8").Insert()

$ws.Range("A8").Value = "QB9BsKic"
$ws.Range("B8").Value = "'01/11/2024"
$ws.Range("C8").Value = "12:15"
$ws.Range("D8").Value = "SAUDI ARABIA - SAUDI PROFESSIONAL LEAGUE"
$ws.Range("E8").Value = "Damac"
$ws.Range("F8").Value = "Al Riyadh"
$ws.Range("G8").Value = 1.9
$ws.Range("H8").Value = 3.4
$ws.Range("I8").Value = 4.2
$ws.Range("J8").Value = 2.4
$ws.Range("K8").Value = 2.2
$ws.Range("L8").Value = 4.33
$ws.Range("M8").Value = 1.05
$ws.Range("N8").Value = 8.5
$ws.Range("O8").Value = 1.25
$ws.Range("P8").Value = 3.75
$ws.Range("Q8").Value = 1.85
$ws.Range("R8").Value = 1.95
$ws.Range("S8").Value = 1.36
$ws.Range("T8").Value = 3
$ws.Range("U8").Value = 1.73
$ws.Range("V8").Value = 2
$ws.Range("W8").Value = 8
$ws.Range("X8").Value = 9
$ws.Range("Y8").Value = 9
$ws.Range("Z8").Value = 15
$ws.Range("AA8").Value = 15
$ws.Range("AB8").Value = 26
$ws.Range("AC8").Value = 10
$ws.Range("AD8").Value = 6.5
$ws.Range("AE8").Value = 13
$ws.Range("AF8").Value = 41
$ws.Range("AG8").Value = 500
$ws.Range("AH8").Value = 13
$ws.Range("AI8").Value = 21
$ws.Range("AJ8").Value = 15
$ws.Range("AK8").Value = 41
$ws.Range("AL8").Value = 34
$ws.Range("AM8").Value = 41
$ws.Range("AN8").Value = 4
$ws.Range("AO8").Value = 10
$ws.Range("AP8").Value = 21
$ws.Range("AQ8").Value = 34
$ws.Range("AR8").Value = 51
$ws.Range("AS8").Value = 126
$ws.Range("AT8").Value = 3
$ws.Range("AU8").Value = 8
$ws.Range("AV8").Value = 51
$ws.Range("AW8").Value = 6
$ws.Range("AX8").Value = 21
$ws.Range("AY8").Value = 29
$ws.Range("AZ8").Value = 67
$ws.Range("BA8").Value = 81
$ws.Range("BB8").Value = 350
$ws.Range("BC8").Value = 81
$ws.Range("BD8").Value = 81
